# Update the dSF (F) column with freshly repulled data.
# (Source data repull: dSF values diverge from dS0 for several outings;
#  this just pushes the newly pulled values into column F.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    11 = 0
    12 = 0
    16 = -1
    21 = -2
    24 = -2
    29 = -1
    31 = 0
    36 = -1
    37 = 0
    54 = -1
    58 = 1
    61 = 1
    67 = -7
    74 = -1
    76 = 1
    78 = 3
    84 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
